$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''57.982.62'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.62%  '

$ws.Range("D3").Value = '''2.282.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.23%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''533.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.14%  '

$ws.Range("D6").Value = '''130.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.50%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  +3.67%  '

$ws.Range("D9").Value = '''2.280.86'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").Value = '''0.0994'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.53%  '

$ws.Range("E11").Value = '  +0.00%  '

$ws.Range("E12").Value = '  +0.76%  '

$ws.Range("D13").Value = '''0.332'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.22%  '

$ws.Range("D14").Value = '''23.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.59%  '

$ws.Range("D15").Value = '''2.688.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.22%  '

$ws.Range("D16").Value = '''57.886.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.75%  '

$ws.Range("E17").Value = '  -0.33%  '

$ws.Range("D18").Value = '''2.285.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("E19").Value = '  -0.87%  '

$ws.Range("E20").Value = '  -2.52%  '

$ws.Range("D21").Value = '''312.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.43%  '

$ws.Range("E22").Value = '  -0.64%  '

$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("E25").Value = '  -0.26%  '

$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("E27").Value = '  -1.64%  '

$ws.Range("E28").Value = '  -2.53%  '

$ws.Range("D29").Value = '''170.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = '''1.70'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.93%  '

$ws.Range("D31").Value = '''0.0' + [char]0x2083 + '0721'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.78%  '

$ws.Range("E32").Value = '  -0.79%  '

$ws.Range("E33").Value = '  -2.02%  '

$ws.Range("E34").Value = '  -0.46%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").Value = '''17.75'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.22%  '

$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("E38").Value = '  -1.05%  '

$ws.Range("E39").Value = '  -1.08%  '

$ws.Range("D40").Value = '''1.49'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.31%  '

$ws.Range("D41").Value = '''140.49'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.38%  '

$ws.Range("D42").Value = '''286.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.06%  '

$ws.Range("D43").Value = '''3.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.22%  '

$ws.Range("D44").Value = '''0.0951'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.00%  '

$ws.Range("D45").Value = '''0.0492'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.55%  '

$ws.Range("E46").Value = '  +0.15%  '

$ws.Range("D47").Value = '''17.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.95%  '

$ws.Range("E48").Value = '  -1.36%  '

$ws.Range("D49").Value = '''10.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.77%  '

$ws.Range("E50").Value = '  -0.48%  '

$ws.Range("E51").Value = '  +1.57%  '
